$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Bump the "Fecha de ingreso" date (column I) by one day for every
#    already-filled visitor row (4-7).
# ---------------------------------------------------------------------
$ws.Range("I4").Value = 46010
$ws.Range("I5").Value = 46010
$ws.Range("I6").Value = 46010
$ws.Range("I7").Value = 46010

# ---------------------------------------------------------------------
# 2. Fill in the three new visitor rows (8, 9 and 10) that were
#    previously blank placeholders, and give the "Nombres" / "RUT" /
#    "Patente Vehicular" cells the smaller 10pt Candara look (no
#    borders) that was used for these freshly-typed-in visitors, while
#    "Empresa", "Cargo", "Motivo de visita" and "Fecha de ingreso" keep
#    reusing the same formatting as the rows above them.
# ---------------------------------------------------------------------

# -- Row 8 : Patricio Sanchez -------------------------------------------------
$ws.Range("B8").Value = "Patricio "
$ws.Range("B8").Font.Name = "Candara"
$ws.Range("B8").Font.Size = 10
$ws.Range("B8").Font.ColorIndex = 1
$ws.Range("B8").Borders.LineStyle = -4142

$ws.Range("C8").Value = "Sanchez "
$ws.Range("E4").Copy()
$ws.Range("C8").PasteSpecial(-4122)

$ws.Range("D8").Value = "15727383-3"
$ws.Range("B8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("E8").Value = "TEK"
$ws.Range("E4").Copy()
$ws.Range("E8").PasteSpecial(-4122)

$ws.Range("F8").Value = "Chofer"
$ws.Range("F4").Copy()
$ws.Range("F8").PasteSpecial(-4122)

$ws.Range("G8").Value = "KRGC-47"
$ws.Range("B8").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$ws.Range("H8").Value = "Ingreso de vehiculos"
$ws.Range("H4").Copy()
$ws.Range("H8").PasteSpecial(-4122)

$ws.Range("I8").Value = 46010
$ws.Range("I4").Copy()
$ws.Range("I8").PasteSpecial(-4122)

# -- Row 9 : Osvaldo Carrasco -------------------------------------------------
$ws.Range("B9").Value = "Osvaldo "
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C9").Value = "Carrasco "
$ws.Range("B8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("D9").Value = "14010414-0 "
$ws.Range("B8").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws.Range("E9").Value = "TEK"
$ws.Range("E4").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("F9").Value = "Chofer"
$ws.Range("F4").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$ws.Range("G9").Value = "RKRD-61"
$ws.Range("B8").Copy()
$ws.Range("G9").PasteSpecial(-4122)

$ws.Range("H9").Value = "Ingreso de vehiculos"
$ws.Range("H4").Copy()
$ws.Range("H9").PasteSpecial(-4122)

$ws.Range("I9").Value = 46010
$ws.Range("I4").Copy()
$ws.Range("I9").PasteSpecial(-4122)

# -- Row 10 : Gonzalo Quezada -------------------------------------------------
$ws.Range("B10").Value = "Gonzalo "
$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("C10").Value = "Quezada "
$ws.Range("B8").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("D10").Value = "12130273-K "
$ws.Range("B8").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("E10").Value = "TEK"
$ws.Range("E4").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("F10").Value = "Chofer"
$ws.Range("F4").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$ws.Range("G10").Value = "KRGC-48"
$ws.Range("B8").Copy()
$ws.Range("G10").PasteSpecial(-4122)

$ws.Range("H10").Value = "Ingreso de vehiculos"
$ws.Range("H4").Copy()
$ws.Range("H10").PasteSpecial(-4122)

$ws.Range("I10").Value = 46010
$ws.Range("I4").Copy()
$ws.Range("I10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Update the active selection to mirror the workbook's last saved
#    state.
# ---------------------------------------------------------------------
$ws.Range("I5:I10").Select()
